$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank data rows 75-77 (Year, Extent on land,
# Extent at sea, Total extent) on the Data_Extent worksheet.
$ws.Range("A75").Value = 2020
$ws.Range("B75").Value = 6.7930987426999998
$ws.Range("C75").Value = 33.821657997800003
$ws.Range("D75").Value = 40.614756740499999

$ws.Range("A76").Value = 2021
$ws.Range("B76").Value = 6.7916245691999988
$ws.Range("C76").Value = 33.823470902699995
$ws.Range("D76").Value = 40.615095471899991

$ws.Range("A77").Value = 2022
$ws.Range("B77").Value = 6.7922804960999992
$ws.Range("C77").Value = 33.823470902699995
$ws.Range("D77").Value = 40.615751398799993

# Reset the view: scroll back to the top-left (clears the stale
# topLeftCell="A58") and move the selection to A2 instead of the old
# D77 / A75:D77 selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
